# Update NATMI LR-pair results with newly computed TPM values.
# A new "ECs" target-cluster row is inserted as the second data row
# (pushing the existing "FAPs" and "MuSCs" target rows down by one),
# and the derived-specificity columns (O, P, S, T) on every row are
# refreshed to reflect the new totals across all target clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the two existing data rows before we start overwriting
# anything, so we can shift them down by one row.
$oldRow2 = @()
for ($c = 1; $c -le 20; $c++) { $oldRow2 += ,$ws.Cells.Item(2, $c).Value() }
$oldRow3 = @()
for ($c = 1; $c -le 20; $c++) { $oldRow3 += ,$ws.Cells.Item(3, $c).Value() }

# Shift: old row 3 ("MuSCs" target) -> row 4, old row 2 ("FAPs" target) -> row 3.
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(4, $c).Value = $oldRow3[$c - 1] }
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(3, $c).Value = $oldRow2[$c - 1] }

# New row 2: Cntn1 -> Ptprz1, FAPs sending cluster -> ECs target cluster.
$ws.Cells.Item(2, 1).Value  = "FAPs"
$ws.Cells.Item(2, 2).Value  = "Cntn1"
$ws.Cells.Item(2, 3).Value  = "Ptprz1"
$ws.Cells.Item(2, 4).Value  = "ECs"
$ws.Cells.Item(2, 5).Value  = 2
$ws.Cells.Item(2, 6).Value  = 0.6666666666666666
$ws.Cells.Item(2, 7).Value  = 0.293933
$ws.Cells.Item(2, 8).Value  = 0.881799
$ws.Cells.Item(2, 9).Value  = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.06163433333333333
$ws.Cells.Item(2, 14).Value = 0.184903
$ws.Cells.Item(2, 15).Value = 0.006690894379667537
$ws.Cells.Item(2, 16).Value = 0.006690894379667537
$ws.Cells.Item(2, 17).Value = 0.01811636449966667
$ws.Cells.Item(2, 18).Value = 0.163047280497
$ws.Cells.Item(2, 19).Value = 0.006690894379667537
$ws.Cells.Item(2, 20).Value = 0.006690894379667537

# Row 3 (the shifted-down "FAPs" target row): only the derived
# specificity columns move, because the denominator now also
# includes the new ECs row.
$ws.Cells.Item(3, 15).Value = 0.01118975126488057
$ws.Cells.Item(3, 16).Value = 0.01118975126488057
$ws.Cells.Item(3, 19).Value = 0.01118975126488057
$ws.Cells.Item(3, 20).Value = 0.01118975126488057

# Row 4 (the shifted-down "MuSCs" target row): same story.
$ws.Cells.Item(4, 15).Value = 0.9821193543554519
$ws.Cells.Item(4, 16).Value = 0.9821193543554518
$ws.Cells.Item(4, 19).Value = 0.9821193543554519
$ws.Cells.Item(4, 20).Value = 0.9821193543554518
